$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.644.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.11%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'3.008.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.59%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.38%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'506.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.29%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'139.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.42%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.35%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.434"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +2.62%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'7.12"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.84%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.108"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.08%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.372"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +5.17%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'3.548.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.11%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.125"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.80%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'26.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +3.57%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.0000163"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +4.93%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'56.652.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.23%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("B17").Value = "'WrappedEther"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'3.028.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.23%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("B18").Value = "'Polkadot"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'6.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.26%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'13.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +4.95%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'7.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.99%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'329.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.98%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.13%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.499"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +3.70%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'65.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +3.87%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'3.157.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.38%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  +0.33%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  +2.06%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'0.0₃0890"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.49%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'6.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.46%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'7.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.43%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'1.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +2.24%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'1.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.59%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'20.48"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.37%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("B34").Value = "'Monero"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'153.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.17%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("B35").Value = "'NEARProtocol"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'4.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.21%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'5.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +2.98%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("B37").Value = "'EnergySwap"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'25.30"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +6.66%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("B38").Value = "'ImmutableX"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'1.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.99%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.0665"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.91%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'3.055.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.14%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'36.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.56%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Value = "'FirstDigitalUSD"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.38%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("B43").Value = "'Filecoin"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'3.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +3.47%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.659"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +3.66%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'2.183.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.19%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'1.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.00%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.946"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.04%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'5.94"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.02%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.0241"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.72%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'19.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +5.65%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").Value = "'Stellar"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'0.0862"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.10%  "
$ws.Range("E51").Style = "Normal"
